# Generate Report for Handoff
# Updates the localization-status report: the handoff file's GUID changed
# from 41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f to 7dd2cd35-1fc4-416b-aab9-5cc4d7d45a55,
# a fresh handoff xliff pair was generated (new hash + timestamps), and the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns reset since no handback has happened yet for the new handoff.

$wb = $excel.ActiveWorkbook

$oldGuid = "41c0a2fa-7ee3-4a4e-8dae-a6c0d4ac5d8f"
$newGuid = "7dd2cd35-1fc4-416b-aab9-5cc4d7d45a55"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newGuid + ".md"
$wsOverview.Range("B2").Value = "e2e\" + $newGuid + ".md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a40e2c853ad572fcfc300acd0f6da4de29d4b9c/e2e/" + $newGuid + ".md", "", "", "e2e\" + $newGuid + ".md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newGuid + ".md"
$wsZh.Range("G2").Value = $newGuid + ".68ca6b393feff937969b3a7eb84ff15764b53f15.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-02 11:10:35"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").ClearFormats()
$wsZh.Range("J2").Value = ""

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a40e2c853ad572fcfc300acd0f6da4de29d4b9c/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md")

$wsZh.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZh.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newGuid + ".md"
$wsDe.Range("G2").Value = $newGuid + ".68ca6b393feff937969b3a7eb84ff15764b53f15.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-02 11:10:39"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").ClearFormats()
$wsDe.Range("J2").Value = ""

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a40e2c853ad572fcfc300acd0f6da4de29d4b9c/e2e/" + $newGuid + ".md", "", "", $newGuid + ".md")

$wsDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDe.Columns.Item(10).ColumnWidth = 20.833333333333332
